# Selenium grid setup: append a second grid-node identity/password pair
# to the credentials sheet and leave the cursor parked past the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$ws.Range("A3").Value = "sumitIdentity2"
$ws.Range("B3").Value = "SummitPass2"

$ws.Range("F4").Select()
